$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7415746450424194
$ws.Range("B1").Value = 1.586796998977661
$ws.Range("C1").Value = 4.804053783416748
$ws.Range("D1").Value = 2.376134634017944
$ws.Range("E1").Value = 1.184810876846313
